# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and refresh
# the associated handoff/generate timestamps on all three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status cells + latest generate date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-18 22:38:54"

# zh-cn sheet: status + latest handoff datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-18 22:38:48"

# de-de sheet: status + latest handoff datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-18 22:38:54"

# Column widths widened to fit the new "Ready for handoff" text
# (ColumnWidth is quantized internally; 16.33 is the input that lands on
# the closest representable width to the target 17.2159881591797)
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
